# TC07_CDS_phs001524_LibrarySource_Genomic.xlsx - "Added CDS All studies testcase"
#
# The SamplesTab query (B3) is trimmed down to drop the Tumor / Analyte Type
# columns (it now only selects Sample ID, Participant ID, Study Name and
# Accession). Re-setting B3's text causes Excel to drop the now-unused old
# string and the FilesTab query (B4, unchanged) naturally slides into its
# place in the shared-strings table - matching the reordering in the diff.
#
# The active selection also moved from B2 to C3 (row 2 is scrolled out of
# view with C3 now the active cell), so we re-select C3 on SamplesTab.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newSamplesTabQuery = @"
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
    s.phs_accession = 'phs001524' AND gi.library_source = 'Genomic'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
"@

# SamplesTab lives in row 3 (A3 = "SamplesTab"); its StatQuery cell is B3.
$ws.Range("B3").Value = $newSamplesTabQuery

# Scroll/selection moved to C3 on the same (Samples) row.
$ws.Range("C3").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
